# Auto-generated cell updates reproducing the cryptos.xlsx refresh commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '70.167.19'
$ws.Range("E2").Value = '  +1.34%  '
# Row 3
$ws.Range("D3").Value = '3.505.36'
$ws.Range("E3").Value = '  +0.21%  '
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.03%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.54'
$ws.Range("E5").Value = '  +0.15%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.18'
$ws.Range("E6").Value = '  +3.13%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.610'
$ws.Range("E7").Value = '  -0.70%  '
# Row 8
$ws.Range("D8").Value = '3.500.04'
$ws.Range("E8").Value = '  +0.18%  '
# Row 9
$ws.Range("E9").Value = '  -0.02%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.192'
$ws.Range("E10").Value = '  -0.24%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.24'
$ws.Range("E11").Value = '  +9.23%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.582'
$ws.Range("E12").Value = '  +0.57%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '46.11'
$ws.Range("E13").Value = '  -1.64%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000274'
$ws.Range("E14").Value = '  -0.97%  '
# Row 15
$ws.Range("D15").Value = '4.069.82'
$ws.Range("E15").Value = '  +0.41%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.29'
$ws.Range("E16").Value = '  +0.01%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '609.80'
$ws.Range("E17").Value = '  -0.05%  '
# Row 18
$ws.Range("D18").Value = '3.506.88'
$ws.Range("E18").Value = '  +0.03%  '
# Row 19
$ws.Range("D19").Value = '70.319.01'
$ws.Range("E19").Value = '  +1.41%  '
# Row 20
$ws.Range("E20").Value = '  +0.81%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.30'
$ws.Range("E21").Value = '  +0.77%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.873'
$ws.Range("E22").Value = '  -0.28%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.02'
$ws.Range("E23").Value = '  -16.28%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.54'
$ws.Range("E24").Value = '  -0.91%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '97.41'
$ws.Range("E25").Value = '  +1.85%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.70'
$ws.Range("E26").Value = '  -3.34%  '
# Row 27
$ws.Range("E27").Value = '  -0.04%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.54'
$ws.Range("E28").Value = '  -1.97%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.70'
$ws.Range("E29").Value = '  +2.40%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.91'
$ws.Range("E30").Value = '  -3.20%  '
# Row 31
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.03'
$ws.Range("E31").Value = '  -4.42%  '
# Row 32
$ws.Range("B32").Value = 'Stacks'
$ws.Range("C32").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.96'
$ws.Range("E32").Value = '  -3.90%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '641.03'
$ws.Range("E33").Value = '  +14.79%  '
# Row 34
$ws.Range("B34").Value = 'Mantle'
$ws.Range("C34").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.27'
$ws.Range("E34").Value = '  -4.45%  '
# Row 35
$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.83'
$ws.Range("E35").Value = '  -0.31%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.56'
$ws.Range("E36").Value = '  +2.28%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0990'
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '10.69'
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0471'
$ws.Range("E39").Value = '  +5.77%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '56.63'
$ws.Range("E40").Value = '  -0.17%  '
# Row 41
$ws.Range("E41").Value = '  +0.03%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.141'
$ws.Range("E42").Value = '  +0.99%  '
# Row 43
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '3.349.28'
$ws.Range("E43").Value = '  -0.38%  '
# Row 44
$ws.Range("B44").Value = 'PEPE'
$ws.Range("C44").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D44").Value = '0.0₃0734'
$ws.Range("E44").Value = '  +5.70%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.307'
$ws.Range("E45").Value = '  -5.29%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '32.13'
$ws.Range("E46").Value = '  -2.01%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.88'
$ws.Range("E47").Value = '  +0.78%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.54'
$ws.Range("E48").Value = '  -2.35%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.129'
$ws.Range("E49").Value = '  +0.62%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.77'
$ws.Range("E50").Value = '  -0.46%  '
